$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values
$ws.Range("B3").Value = 73758936550105.94
$ws.Range("C3").Value = 54669803495241.19
$ws.Range("D3").Value = 699264329849554.4

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.04225674932233511
$ws.Range("C4").Value = 0.04101156979433512
$ws.Range("D4").Value = 292692101483443.2

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 130444170256473.6
$ws.Range("C5").Value = 23853589868413.27
$ws.Range("D5").Value = 253807614762507.3
